$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing values (type_lieu / Agence column D emptied, other fields changed) ---
$ws.Range("A2").Value = "mmmmml"

$ws.Range("B2").Value = "'101"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "'011111111000101010100101"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "'"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "'10"
$ws.Range("E2").Style = "Normal"

$ws.Range("G2").Value = "444/AA4444"

$ws.Range("I2").Value = 12200
$ws.Range("J2").Value = 1830
$ws.Range("K2").Value = 10370

# --- Row 3: brand new row appended below ---
$ws.Range("A3").Value = "aaaaaaa"
$ws.Range("B3").Value = "aaaaa"

$ws.Range("C3").Value = "'121313213213213213213213"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "'321321"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "006/tESTDRR"
$ws.Range("H3").Value = "mensuelle"

$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 8500
